# Loan RBI, Variable Instalments
# - Insert a new (blank) column before column N on the "Repayment Schedule"
#   sheet, pushing the existing "Late" / "Outstanding" columns one to the
#   right (N->O, P->Q), and give the new column the same width as column M.
# - Make "Repayment Schedule" the active sheet/tab (was "Transactions"),
#   with a new selection of R9.

$wb = $excel.ActiveWorkbook

$repayment = $wb.Worksheets.Item("Repayment Schedule")

# Insert a blank column in front of column N (shifts N:P -> O:Q).
$repayment.Columns("N:N").Insert()

# The newly inserted column picks up the width of the column to its left.
$repayment.Columns("N:N").ColumnWidth = $repayment.Columns("M:M").ColumnWidth

# Switch the active sheet to "Repayment Schedule" (was "Transactions") and
# move the selection to R9.
$repayment.Activate()
[void]$repayment.Range("R9").Select()
